# ETHICS grades quiz 3
# Update quiz grades (column F) for assignment #2 and #3 on Sheet1,
# then move the active selection to J13 to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Assignment #2 (row 9) quiz grade: 1 -> 0.9
$ws.Range("F9").Value = 0.9

# Assignment #3 (row 10) quiz grade: 1 -> 0.7
$ws.Range("F10").Value = 0.7

# Update the selection to match the recorded cursor position
$ws.Activate()
$ws.Range("J13").Select()
